$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Activate()
$summary.Range("B4").Value = 0
$summary.Range("E4").Value = 200

# --- Repayment schedule sheet updates ---
$repay = $wb.Worksheets.Item("Repayment schedule")
$repay.Activate()

$repay.Range("P2").Clear()

$repay.Range("K3").Value = 937.72
$repay.Range("O3").Clear()
$repay.Range("P3").Value = 937.72

$repay.Range("O4").Clear()

$repay.Range("K5").Value = 937.72
$repay.Range("O5").Clear()
$repay.Range("P5").Value = 937.72

$repay.Range("O6").Clear()

$repay.Range("O7").Clear()

$repay.Range("O8").Clear()

# Final selections to match the recorded view state:
# Summary sheet cursor on F4 (no longer the active tab)...
$summary.Range("F4").Select()
# ...and Repayment schedule is the active tab with H4 selected.
$repay.Activate()
$repay.Range("H4").Select()
